$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.181047201156616
$ws.Range("B1").Value = 2.190061807632446
$ws.Range("C1").Value = 10.49327087402344
$ws.Range("D1").Value = 2.56791090965271
$ws.Range("E1").Value = 1.241783857345581
